$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# May (row 5) and June (row 6) 2023 circulation figures added to column F
$ws.Range("F5").Value = 84778
$ws.Range("F6").Value = 95894

# Move/leave the active selection on F7, matching the author's cursor position
$ws.Range("F7").Select()
